$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 10: APMS-T153 / Ankit@123 ---
$ws.Range("A10").Value = "APMS-T153"
$ws.Range("C2").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("C10").Value = "Ankit@123"
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:Ankit@123")
$ws.Range("C6").Copy()
$ws.Range("C10").PasteSpecial(-4122)

# --- Row 11: APMS-T154 / AnkitYadav@1234 ---
$ws.Range("A11").Value = "APMS-T154"
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("C11").Value = "AnkitYadav@1234"
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:AnkitYadav@1234")
$ws.Range("C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- Row 12: APMS-T155 / AnkitYadav@12345$ ---
$ws.Range("A12").Value = "APMS-T155"
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("C12").Value = "AnkitYadav@12345$"
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:AnkitYadav@12345$")
$ws.Range("C6").Copy()
$ws.Range("C12").PasteSpecial(-4122)

# Update the active selection to match the saved view state
[void]$ws.Range("D7").Select()
